$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.332.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.98%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''3.326.85'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -4.78%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = '''  -0.03%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = '''572.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.93%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = '''176.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +2.23%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("B7").Value = '''USDC'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = '''https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +0.03%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("B8").Value = '''XRP'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = '''https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = '''0.609'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +1.63%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''3.318.76'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -5.01%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = '''0.128'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -2.25%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = '''  -0.25%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = '''  -1.09%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = '''3.896.31'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -4.92%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = '''  +0.23%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = '''28.38'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -5.25%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = '''65.356.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -1.08%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = '''0.0000167'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -2.04%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''3.333.59'
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = '''  -3.53%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = '''13.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -4.28%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = '''361.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -1.45%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = '''7.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -4.55%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = '''0.997'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -0.38%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = '''71.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -1.68%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = '''0.516'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -3.74%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = '''  -4.14%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = '''9.48'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -1.46%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = '''  -1.15%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = '''  +0.01%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = '''1.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -1.77%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = '''1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.07%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = '''5.57'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -3.64%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = '''22.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -5.28%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = '''6.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -4.71%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = '''  -6.76%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = '''Monero'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = '''160.24'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +0.56%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = '''ImmutableX'
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = '''1.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -3.86%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = '''0.840'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -5.73%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = '''27.22'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -8.14%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = '''  -0.99%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = '''Maker'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = '''2.701.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -4.49%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = '''dogwifhat'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = '''https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = '''2.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -1.82%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = '''6.19'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -4.89%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = '''4.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -4.63%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = '''Hedera'
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = '''0.0665'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -2.57%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = '''OKB'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = '''39.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -1.41%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = '''Bittensor'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''333.45'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +3.61%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = '''23.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -1.37%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = '''0.0277'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -3.89%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = '''0.102'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +1.22%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = '''FirstDigitalUSD'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = '''0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.00%  '
$ws.Range("E51").Style = "Normal"
